$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 58: "1372. Longest ZigZag Path in a Binary Tree" (Trees, Medium)
$question = '1372. Longest ZigZag Path in a Binary Tree'
$difficulty = 'Medium'
$pattern = 'Trees'
$notes = 'For a simple solution, maintain a global max variable. Propagate the depth and the direction. Call the recursive function on both sides at the start, and at each node, check to override max, and then call dfs on both sides once again. If the conditions are maintained, add depth+1, else reset depth to 0.'
$link = 'https://leetcode.com/problems/longest-zigzag-path-in-a-binary-tree/solutions/531867/java-python-dfs-solution/?envType=study-plan-v2&envId=leetcode-75 '

# Column A - Question
$ws.Range("A58").Value = $question

# Column E - Link, styled like the other hyperlink cells (copy format from E57), then add the hyperlink
# (set before D58/Notes so the shared-string table gets new entries in the same order as the source edit)
$ws.Range("E57").Copy()
$ws.Range("E58").PasteSpecial(-4122)
$ws.Range("E58").Value = $link
$ws.Hyperlinks.Add($ws.Range("E58"), $link)
$ws.Range("E57").Copy()
$ws.Range("E58").PasteSpecial(-4122)

# Column D - Notes
$ws.Range("D58").Value = $notes

# Column B - Difficulty, matching the "Medium" highlight style used elsewhere (copy format from B57)
$ws.Range("B57").Copy()
$ws.Range("B58").PasteSpecial(-4122)
$ws.Range("B58").Value = $difficulty

# Column C - Pattern
$ws.Range("C58").Value = $pattern

$excel.CutCopyMode = $false

# Update selection to match the saved view state
$ws.Range("E62").Select() | Out-Null
